# Add a new parameter row ("track_region") to the workflow_config sheet,
# describing the genomic region to plot genome tracks over.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# Insert a new row above the current row 12 ("ensembl_version"), shifting
# every row below it down by one (12-28 -> 13-29).
$ws.Rows.Item(12).Insert()

# Pick up the same visual formatting (borders/fill/font/number format) that
# the rest of the parameter table uses by copying it from the row that used
# to be 12 and is now 13 (style indices only - values are untouched).
$ws.Range("A13:C13").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New row 12 is shorter (2 lines) than the old row 12 (48pt / 3 lines), so
# give it the 32pt height used by the other 2-line rows in this table.
$ws.Rows.Item(12).RowHeight = 32

# --- Cell contents -----------------------------------------------------
$ws.Range("A12").Value = "track_region"

$descText = "Genomic region to plot genome tracks over. Example: chr1:500000-900000 (REQUIRED if run_genome_track)"
$ws.Range("B12").Value = $descText

$ws.Range("C12").Value = "chr1:700000-900000"

# --- Rich text formatting for the description cell ---------------------
# Bold "Example: chr1:500000-900000"
$exampleStart = $descText.IndexOf("Example:") + 1
$exampleLen = "Example: chr1:500000-900000".Length
$ws.Range("B12").Characters($exampleStart, $exampleLen).Font.Bold = $true

# Bold "(REQUIRED if run_genome_track)"
$reqStart = $descText.IndexOf("(REQUIRED if run_genome_track)") + 1
$reqLen = "(REQUIRED if run_genome_track)".Length
$ws.Range("B12").Characters($reqStart, $reqLen).Font.Bold = $true

# Match the post-edit selection recorded in the workbook (user left the
# selection on C13, the cell that used to be C12 before the insert).
$ws.Range("C13").Select()
